$d = $word.ActiveDocument

$pairs = @(
    @("2025-07-17 Thursday", "2025-07-18 Friday"),
    @("397×4=1588", "362×9=3258"),
    @("945×3=2835", "250×8=2000"),
    @("931×9=8379", "639×7=4473"),
    @("148×3=444", "364×6=2184"),
    @("963×6=5778", "815×8=6520"),
    @("672×8=5376", "823×8=6584"),
    @("907×3=2721", "326×8=2608"),
    @("556×5=2780", "297×4=1188"),
    @("317×2=634", "722×8=5776"),
    @("557×6=3342", "788×9=7092"),
    @("792×4=3168", "754×3=2262"),
    @("649×3=1947", "120×6=720"),
    @("498×8=3984", "683×7=4781"),
    @("308×5=1540", "860×3=2580"),
    @("459×7=3213", "985×2=1970"),
    @("749×7=5243", "169×7=1183"),
    @("302×4=1208", "395×8=3160"),
    @("337×8=2696", "765×5=3825"),
    @("248×8=1984", "607×3=1821"),
    @("185×7=1295", "784×7=5488"),
    @("338×4=1352", "180×9=1620"),
    @("342×7=2394", "963×7=6741"),
    @("356×6=2136", "765×2=1530"),
    @("382×7=2674", "464×3=1392"),
    @("829×6=4974", "479×4=1916")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
